# Scheduled market-data refresh: update currentAveragePrice / Leve price & profit
# columns (H:N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 27782422
$ws.Range("I41").Value = 967.8333
$ws.Range("J41").Value = 41673148
$ws.Range("K41").Value = 967.8333
$ws.Range("L41").Value = 41673148
$ws.Range("M41").Value = -527.8333
$ws.Range("N41").Value = -41674028

$ws.Range("H94").Value = 1977.8572
$ws.Range("I94").Value = 1478.3334
$ws.Range("J94").Value = 4975
$ws.Range("K94").Value = 1478.3334
$ws.Range("L94").Value = 4975
$ws.Range("M94").Value = -1027.3334
$ws.Range("N94").Value = -5877

$ws.Range("H96").Value = 675.5
$ws.Range("I96").Value = 639.44446
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1918.33338
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -545.33338
$ws.Range("N96").Value = -5746

$ws.Range("H98").Value = 1236.6957
$ws.Range("I98").Value = 797.25
$ws.Range("J98").Value = 4166.3335
$ws.Range("K98").Value = 797.25
$ws.Range("L98").Value = 4166.3335
$ws.Range("M98").Value = 700.75

$ws.Range("H103").Value = 729.1429000000001
$ws.Range("I103").Value = 868
$ws.Range("J103").Value = 625
$ws.Range("K103").Value = 2604
$ws.Range("L103").Value = 1875
$ws.Range("M103").Value = -2018
$ws.Range("N103").Value = -3047

$ws.Range("H122").Value = 1236.6957
$ws.Range("I122").Value = 797.25
$ws.Range("J122").Value = 4166.3335
$ws.Range("K122").Value = 2391.75
$ws.Range("L122").Value = 12499.0005
$ws.Range("M122").Value = 58.25

$ws.Range("H129").Value = 1688.9166
$ws.Range("I129").Value = 717.61536
$ws.Range("J129").Value = 2836.818
$ws.Range("K129").Value = 2152.84608
$ws.Range("L129").Value = 8510.454000000002
$ws.Range("M129").Value = 2847.15392
$ws.Range("N129").Value = -18510.454

$ws.Range("H132").Value = 2977.6155
$ws.Range("I132").Value = 2858.1
$ws.Range("J132").Value = 3376
$ws.Range("K132").Value = 8574.299999999999
$ws.Range("L132").Value = 10128
$ws.Range("M132").Value = -6044.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3582.5154
$ws.Range("I32").Value = 3593.802
$ws.Range("J32").Value = 2499
$ws.Range("K32").Value = 3593.802
$ws.Range("L32").Value = 2499
$ws.Range("M32").Value = -3306.802

$ws.Range("H97").Value = 1106.2916
$ws.Range("I97").Value = 1093.5217
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 1093.5217
$ws.Range("L97").Value = 1400
$ws.Range("M97").Value = -597.5217
$ws.Range("N97").Value = -2392

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1123.3846
$ws.Range("I94").Value = 1184.4546
$ws.Range("J94").Value = 787.5
$ws.Range("K94").Value = 1184.4546
$ws.Range("L94").Value = 787.5
$ws.Range("M94").Value = -733.4546
$ws.Range("N94").Value = -1689.5

$ws.Range("H105").Value = 3776.3333
$ws.Range("I105").Value = 3776.3333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3776.3333
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -2029.3333

$ws.Range("H107").Value = 627613.2
$ws.Range("I107").Value = 2304
$ws.Range("J107").Value = 3337286.2
$ws.Range("K107").Value = 2304
$ws.Range("L107").Value = 3337286.2
$ws.Range("M107").Value = -384
$ws.Range("N107").Value = -3341126.2

$ws.Range("H134").Value = 56090.684
$ws.Range("I134").Value = 3651.2778
$ws.Range("J134").Value = 1000000
$ws.Range("K134").Value = 10953.8334
$ws.Range("L134").Value = 3000000
$ws.Range("M134").Value = -8418.8334
$ws.Range("N134").Value = -3005070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 102526.9
$ws.Range("I31").Value = 1283.75
$ws.Range("J31").Value = 507499.5
$ws.Range("K31").Value = 1283.75
$ws.Range("L31").Value = 507499.5
$ws.Range("M31").Value = -988.75
$ws.Range("N31").Value = -508089.5

$ws.Range("H34").Value = 102526.9
$ws.Range("I34").Value = 1283.75
$ws.Range("J34").Value = 507499.5
$ws.Range("K34").Value = 1283.75
$ws.Range("L34").Value = 507499.5
$ws.Range("M34").Value = -1081.75
$ws.Range("N34").Value = -507903.5

$ws.Range("H107").Value = 657.88464
$ws.Range("I107").Value = 534.2381
$ws.Range("J107").Value = 1177.2
$ws.Range("K107").Value = 534.2381
$ws.Range("L107").Value = 1177.2
$ws.Range("M107").Value = 1385.7619

$ws.Range("H132").Value = 2405.15
$ws.Range("I132").Value = 2014.5385
$ws.Range("J132").Value = 3130.5715
$ws.Range("K132").Value = 6043.6155
$ws.Range("L132").Value = 9391.7145
$ws.Range("M132").Value = -3513.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 43143.207
$ws.Range("I122").Value = 825.7143
$ws.Range("J122").Value = 60568.06
$ws.Range("K122").Value = 7431.428699999999
$ws.Range("L122").Value = 545112.54
$ws.Range("M122").Value = -4981.428699999999
$ws.Range("N122").Value = -550012.54

$ws.Range("H127").Value = 1622
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 1622
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 4866
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -14786

$ws.Range("H131").Value = 2901.0195
$ws.Range("I131").Value = 1678.091
$ws.Range("J131").Value = 3237.325
$ws.Range("K131").Value = 5034.272999999999
$ws.Range("L131").Value = 9711.974999999999
$ws.Range("M131").Value = 5.727000000000771
$ws.Range("N131").Value = -19791.975

$ws.Range("H133").Value = 19059.613
$ws.Range("I133").Value = 7808.1665
$ws.Range("J133").Value = 21759.96
$ws.Range("K133").Value = 23424.4995
$ws.Range("L133").Value = 65279.88
$ws.Range("M133").Value = -18364.4995

$ws.Range("H138").Value = 4051.7
$ws.Range("I138").Value = 3704
$ws.Range("J138").Value = 4399.4
$ws.Range("K138").Value = 11112
$ws.Range("L138").Value = 13198.2
$ws.Range("M138").Value = -5972
$ws.Range("N138").Value = -23478.2

$ws.Range("H139").Value = 9928.5
$ws.Range("I139").Value = 8999
$ws.Range("J139").Value = 10000
$ws.Range("K139").Value = 26997
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -21857

$ws.Range("H140").Value = 7020
$ws.Range("I140").Value = 4500
$ws.Range("J140").Value = 7200
$ws.Range("K140").Value = 13500
$ws.Range("L140").Value = 21600
$ws.Range("M140").Value = -8320
$ws.Range("N140").Value = -31960

$ws.Range("H141").Value = 4284.143
$ws.Range("I141").Value = 4284.143
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 12852.429
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -7672.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7037.8667
$ws.Range("I70").Value = 5729.0454
$ws.Range("J70").Value = 10637.125
$ws.Range("K70").Value = 5729.0454
$ws.Range("L70").Value = 10637.125
$ws.Range("M70").Value = -5459.0454

$ws.Range("H73").Value = 7037.8667
$ws.Range("I73").Value = 5729.0454
$ws.Range("J73").Value = 10637.125
$ws.Range("K73").Value = 5729.0454
$ws.Range("L73").Value = 10637.125
$ws.Range("M73").Value = -4793.0454

$ws.Range("H113").Value = 359400.16
$ws.Range("I113").Value = 527623.9
$ws.Range("J113").Value = 4261.1113
$ws.Range("K113").Value = 527623.9
$ws.Range("L113").Value = 4261.1113
$ws.Range("M113").Value = -525453.9

$ws.Range("H132").Value = 34469.344
$ws.Range("I132").Value = 3441.4138
$ws.Range("J132").Value = 334406
$ws.Range("K132").Value = 10324.2414
$ws.Range("L132").Value = 1003218
$ws.Range("M132").Value = -7794.241399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1719.8
$ws.Range("I22").Value = 1749.75
$ws.Range("J22").Value = 1600
$ws.Range("K22").Value = 1749.75
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = -1454.75
$ws.Range("N22").Value = -2190

$ws.Range("H27").Value = 1719.8
$ws.Range("I27").Value = 1749.75
$ws.Range("J27").Value = 1600
$ws.Range("K27").Value = 1749.75
$ws.Range("L27").Value = 1600
$ws.Range("M27").Value = -1642.75
$ws.Range("N27").Value = -1814

$ws.Range("H40").Value = 4824.4614
$ws.Range("I40").Value = 4068.3809
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 4068.3809
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -3932.3809

$ws.Range("H46").Value = 2212
$ws.Range("I46").Value = 2620.889
$ws.Range("J46").Value = 1598.6666
$ws.Range("K46").Value = 2620.889
$ws.Range("L46").Value = 1598.6666
$ws.Range("M46").Value = -2432.889
$ws.Range("N46").Value = -1974.6666

$ws.Range("H50").Value = 38040
$ws.Range("I50").Value = 60000
$ws.Range("J50").Value = 33648
$ws.Range("K50").Value = 60000
$ws.Range("L50").Value = 33648
$ws.Range("M50").Value = -59363
$ws.Range("N50").Value = -34922

$ws.Range("H55").Value = 2441
$ws.Range("I55").Value = 79.333336
$ws.Range("J55").Value = 4212.25
$ws.Range("K55").Value = 79.333336
$ws.Range("L55").Value = 4212.25
$ws.Range("M55").Value = 93.666664
$ws.Range("N55").Value = -4558.25

$ws.Range("H61").Value = 2769.2727
$ws.Range("I61").Value = 2946.55
$ws.Range("J61").Value = 996.5
$ws.Range("K61").Value = 2946.55
$ws.Range("L61").Value = 996.5
$ws.Range("M61").Value = -2744.55
$ws.Range("N61").Value = -1400.5

$ws.Range("H82").Value = 1294.8
$ws.Range("I82").Value = 1319
$ws.Range("J82").Value = 1228.25
$ws.Range("K82").Value = 1319
$ws.Range("L82").Value = 1228.25
$ws.Range("M82").Value = -958

$ws.Range("H85").Value = 1294.8
$ws.Range("I85").Value = 1319
$ws.Range("J85").Value = 1228.25
$ws.Range("K85").Value = 1319
$ws.Range("L85").Value = 1228.25
$ws.Range("M85").Value = -71

$ws.Range("H113").Value = 2769.2727
$ws.Range("I113").Value = 2946.55
$ws.Range("J113").Value = 996.5
$ws.Range("K113").Value = 2946.55
$ws.Range("L113").Value = 996.5
$ws.Range("M113").Value = -776.5500000000002
$ws.Range("N113").Value = -5336.5

$ws.Range("H136").Value = 222832.73
$ws.Range("I136").Value = 361112.5
$ws.Range("J136").Value = 7730.8887
$ws.Range("K136").Value = 1083337.5
$ws.Range("L136").Value = 23192.6661
$ws.Range("M136").Value = -1080787.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 895.75
$ws.Range("I100").Value = 810.61536
$ws.Range("J100").Value = 1264.6666
$ws.Range("K100").Value = 1621.23072
$ws.Range("L100").Value = 2529.3332
$ws.Range("M100").Value = -1080.23072
$ws.Range("N100").Value = -3611.3332

$ws.Range("H107").Value = 1432.16
$ws.Range("I107").Value = 1616.25
$ws.Range("J107").Value = 695.8
$ws.Range("K107").Value = 4848.75
$ws.Range("L107").Value = 2087.4
$ws.Range("M107").Value = -2928.75

$ws.Range("H126").Value = 1978
$ws.Range("I126").Value = 1914.5
$ws.Range("J126").Value = 2105
$ws.Range("K126").Value = 5743.5
$ws.Range("L126").Value = 6315
$ws.Range("M126").Value = -3273.5

$ws.Range("H136").Value = 82656.39999999999
$ws.Range("I136").Value = 2720.35
$ws.Range("J136").Value = 402400.6
$ws.Range("K136").Value = 8161.049999999999
$ws.Range("L136").Value = 1207201.8
$ws.Range("M136").Value = -5611.049999999999
$ws.Range("N136").Value = -1212301.8
